$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 328
$ws.Range("A328").Value = 327
$ws.Range("B328").Value = "Friday, Jan 13"
$ws.Range("C328").Value = "12:13 AM"
$ws.Range("D328").Value = "SAR142"
$ws.Range("E328").Value = "Karlsruhe/Baden-Baden"
$ws.Range("F328").Value = "(FKB)"
$ws.Range("G328").Value = "SprintAir "
$ws.Range("H328").Value = "SF34"
$ws.Range("I328").Value = "(SP-KPE)"
$ws.Range("J328").Value = "1:12 AM"
$ws.Range("L328").Value = "0 hours, 59 minutes"

# Row 329
$ws.Range("A329").Value = 328
$ws.Range("B329").Value = "Friday, Jan 13"
$ws.Range("C329").Value = "5:48 AM"
$ws.Range("D329").Value = "UNKNOWN"
$ws.Range("E329").Value = "Leipzig"
$ws.Range("F329").Value = "(LEJ)"
$ws.Range("G329").Value = "DHL "
$ws.Range("H329").Value = "B738"
$ws.Range("I329").Value = "(EC-IXO)"
$ws.Range("J329").Value = "5:51 AM"
$ws.Range("L329").Value = "0 hours, 3 minutes"

# Row 330
$ws.Range("A330").Value = 329
$ws.Range("B330").Value = "Friday, Jan 13"
$ws.Range("C330").Value = "7:35 AM"
$ws.Range("D330").Value = "FR4105"
$ws.Range("E330").Value = "Wroclaw"
$ws.Range("F330").Value = "(WRO)"
$ws.Range("G330").Value = "Ryanair "
$ws.Range("H330").Value = "B738"
$ws.Range("I330").Value = "(SP-RKG)"
$ws.Range("J330").Value = "7:23 AM"
$ws.Range("L330").Value = "0 hours, -12 minutes"

# Row 331
$ws.Range("A331").Value = 330
$ws.Range("B331").Value = "Friday, Jan 13"
$ws.Range("C331").Value = "7:55 AM"
$ws.Range("D331").Value = "3V4006"
$ws.Range("E331").Value = "Prague"
$ws.Range("F331").Value = "(PRG)"
$ws.Range("G331").Value = "ASL Airlines "
$ws.Range("H331").Value = "B738"
$ws.Range("I331").Value = "(OE-IXH)"
$ws.Range("J331").Value = "8:01 AM"
$ws.Range("L331").Value = "0 hours, 6 minutes"

# Row 332
$ws.Range("A332").Value = 331
$ws.Range("B332").Value = "Friday, Jan 13"
$ws.Range("C332").Value = "8:08 AM"
$ws.Range("D332").Value = "P81956"
$ws.Range("E332").Value = "Berlin"
$ws.Range("F332").Value = "(BER)"
$ws.Range("G332").Value = "SprintAir "
$ws.Range("H332").Value = "SF34"
$ws.Range("I332").Value = "(SP-KPV)"
$ws.Range("J332").Value = "8:18 AM"
$ws.Range("L332").Value = "0 hours, 10 minutes"

# Row 333
$ws.Range("A333").Value = 332
$ws.Range("B333").Value = "Friday, Jan 13"
$ws.Range("C333").Value = "8:10 AM"
$ws.Range("D333").Value = "LO3831"
$ws.Range("E333").Value = "Warsaw"
$ws.Range("F333").Value = "(WAW)"
$ws.Range("G333").Value = "LOT (Warmia Mazury Livery) "
$ws.Range("H333").Value = "E195"
$ws.Range("I333").Value = "(SP-LNF)"
$ws.Range("J333").Value = "7:56 AM"
$ws.Range("L333").Value = "0 hours, -14 minutes"

# Row 334
$ws.Range("A334").Value = 333
$ws.Range("B334").Value = "Friday, Jan 13"
$ws.Range("C334").Value = "9:05 AM"
$ws.Range("D334").Value = "FR6120"
$ws.Range("E334").Value = "Hamburg"
$ws.Range("F334").Value = "(HAM)"
$ws.Range("G334").Value = "Ryanair "
$ws.Range("H334").Value = "B738"
$ws.Range("I334").Value = "(SP-RSW)"
$ws.Range("J334").Value = "9:08 AM"
$ws.Range("L334").Value = "0 hours, 3 minutes"

# Row 335
$ws.Range("A335").Value = 334
$ws.Range("B335").Value = "Friday, Jan 13"
$ws.Range("C335").Value = "9:10 AM"
$ws.Range("D335").Value = "FR6085"
$ws.Range("E335").Value = "London"
$ws.Range("F335").Value = "(STN)"
$ws.Range("G335").Value = "Ryanair "
$ws.Range("H335").Value = "B38M"
$ws.Range("I335").Value = "(EI-HGT)"
$ws.Range("J335").Value = "8:58 AM"
$ws.Range("L335").Value = "0 hours, -12 minutes"

# Row 336
$ws.Range("A336").Value = 335
$ws.Range("B336").Value = "Friday, Jan 13"
$ws.Range("C336").Value = "9:15 AM"
$ws.Range("D336").Value = "FR6845"
$ws.Range("E336").Value = "Copenhagen"
$ws.Range("F336").Value = "(CPH)"
$ws.Range("G336").Value = "Ryanair "
$ws.Range("H336").Value = "B738"
$ws.Range("I336").Value = "(SP-RKQ)"
$ws.Range("J336").Value = "9:04 AM"
$ws.Range("L336").Value = "0 hours, -11 minutes"

# Row 337
$ws.Range("A337").Value = 336
$ws.Range("B337").Value = "Friday, Jan 13"
$ws.Range("C337").Value = "9:20 AM"
$ws.Range("D337").Value = "FR6102"
$ws.Range("E337").Value = "Stockholm"
$ws.Range("F337").Value = "(ARN)"
$ws.Range("G337").Value = "Ryanair "
$ws.Range("H337").Value = "B738"
$ws.Range("I337").Value = "(SP-RKM)"
$ws.Range("J337").Value = "9:19 AM"
$ws.Range("L337").Value = "0 hours, -1 minutes"

# Row 338
$ws.Range("A338").Value = 337
$ws.Range("B338").Value = "Friday, Jan 13"
$ws.Range("C338").Value = "9:30 AM"
$ws.Range("D338").Value = "FR6098"
$ws.Range("E338").Value = "Gothenburg"
$ws.Range("F338").Value = "(GOT)"
$ws.Range("G338").Value = "Ryanair "
$ws.Range("H338").Value = "B738"
$ws.Range("I338").Value = "(SP-RSO)"
$ws.Range("J338").Value = "9:15 AM"
$ws.Range("L338").Value = "0 hours, -15 minutes"

# Row 339
$ws.Range("A339").Value = 338
$ws.Range("B339").Value = "Friday, Jan 13"
$ws.Range("C339").Value = "10:35 AM"
$ws.Range("D339").Value = "W61642"
$ws.Range("E339").Value = "Eindhoven"
$ws.Range("F339").Value = "(EIN)"
$ws.Range("G339").Value = "Wizz Air "
$ws.Range("H339").Value = "A21N"
$ws.Range("I339").Value = "(9H-WBU)"
$ws.Range("J339").Value = "10:33 AM"
$ws.Range("L339").Value = "0 hours, -2 minutes"

# Row 340
$ws.Range("A340").Value = 339
$ws.Range("B340").Value = "Friday, Jan 13"
$ws.Range("C340").Value = "11:25 AM"
$ws.Range("D340").Value = "LO3835"
$ws.Range("E340").Value = "Warsaw"
$ws.Range("F340").Value = "(WAW)"
$ws.Range("G340").Value = "LOT "
$ws.Range("H340").Value = "E170"
$ws.Range("I340").Value = "(SP-LDH)"
$ws.Range("J340").Value = "11:11 AM"
$ws.Range("L340").Value = "0 hours, -14 minutes"

# Row 341
$ws.Range("A341").Value = 340
$ws.Range("B341").Value = "Friday, Jan 13"
$ws.Range("C341").Value = "11:30 AM"
$ws.Range("D341").Value = "W61602"
$ws.Range("E341").Value = "London"
$ws.Range("F341").Value = "(LTN)"
$ws.Range("G341").Value = "Wizz Air "
$ws.Range("H341").Value = "A320"
$ws.Range("I341").Value = "(HA-LYS)"
$ws.Range("J341").Value = "11:31 AM"
$ws.Range("L341").Value = "0 hours, 1 minutes"

# Row 342
$ws.Range("A342").Value = 341
$ws.Range("B342").Value = "Friday, Jan 13"
$ws.Range("C342").Value = "12:20 PM"
$ws.Range("D342").Value = "DY1030"
$ws.Range("E342").Value = "Bergen"
$ws.Range("F342").Value = "(BGO)"
$ws.Range("G342").Value = "Norwegian "
$ws.Range("H342").Value = "B738"
$ws.Range("I342").Value = "(LN-ENM)"
$ws.Range("J342").Value = "12:07 PM"
$ws.Range("L342").Value = "0 hours, -13 minutes"

# Row 343
$ws.Range("A343").Value = 342
$ws.Range("B343").Value = "Friday, Jan 13"
$ws.Range("C343").Value = "12:20 PM"
$ws.Range("D343").Value = "FR8781"
$ws.Range("E343").Value = "Cork"
$ws.Range("F343").Value = "(ORK)"
$ws.Range("G343").Value = "Ryanair "
$ws.Range("H343").Value = "B738"
$ws.Range("I343").Value = "(SP-RSL)"
$ws.Range("J343").Value = "11:57 AM"
$ws.Range("L343").Value = "0 hours, -23 minutes"

# Row 344
$ws.Range("A344").Value = 343
$ws.Range("B344").Value = "Friday, Jan 13"
$ws.Range("C344").Value = "12:30 PM"
$ws.Range("D344").Value = "DY1052"
$ws.Range("E344").Value = "Oslo"
$ws.Range("F344").Value = "(OSL)"
$ws.Range("G344").Value = "Norwegian "
$ws.Range("H344").Value = "B738"
$ws.Range("I344").Value = "(SE-RPG)"
$ws.Range("J344").Value = "12:25 PM"
$ws.Range("L344").Value = "0 hours, -5 minutes"

# Row 345
$ws.Range("A345").Value = 344
$ws.Range("B345").Value = "Friday, Jan 13"
$ws.Range("C345").Value = "12:35 PM"
$ws.Range("D345").Value = "FR3687"
$ws.Range("E345").Value = "Billund"
$ws.Range("F345").Value = "(BLL)"
$ws.Range("G345").Value = "Ryanair "
$ws.Range("H345").Value = "B738"
$ws.Range("I345").Value = "(SP-RSW)"
$ws.Range("J345").Value = "12:38 PM"
$ws.Range("L345").Value = "0 hours, 3 minutes"

# Row 346
$ws.Range("A346").Value = 345
$ws.Range("B346").Value = "Friday, Jan 13"
$ws.Range("C346").Value = "12:35 PM"
$ws.Range("D346").Value = "LH1642"
$ws.Range("E346").Value = "Munich"
$ws.Range("F346").Value = "(MUC)"
$ws.Range("G346").Value = "Lufthansa "
$ws.Range("H346").Value = "CRJ9"
$ws.Range("I346").Value = "(D-ACNX)"
$ws.Range("J346").Value = "12:21 PM"
$ws.Range("L346").Value = "0 hours, -14 minutes"

# Row 347
$ws.Range("A347").Value = 346
$ws.Range("B347").Value = "Friday, Jan 13"
$ws.Range("C347").Value = "1:05 PM"
$ws.Range("D347").Value = "FR3278"
$ws.Range("E347").Value = "Oslo"
$ws.Range("F347").Value = "(TRF)"
$ws.Range("G347").Value = "Ryanair "
$ws.Range("H347").Value = "B738"
$ws.Range("I347").Value = "(SP-RKQ)"
$ws.Range("J347").Value = "1:15 PM"
$ws.Range("L347").Value = "0 hours, 10 minutes"

# Row 348
$ws.Range("A348").Value = 347
$ws.Range("B348").Value = "Friday, Jan 13"
$ws.Range("C348").Value = "1:30 PM"
$ws.Range("D348").Value = "LH1376"
$ws.Range("E348").Value = "Frankfurt"
$ws.Range("F348").Value = "(FRA)"
$ws.Range("G348").Value = "Lufthansa "
$ws.Range("H348").Value = "CRJ9"
$ws.Range("I348").Value = "(D-ACNR)"
$ws.Range("J348").Value = "1:25 PM"
$ws.Range("L348").Value = "0 hours, -5 minutes"

# Row 349
$ws.Range("A349").Value = 348
$ws.Range("B349").Value = "Friday, Jan 13"
$ws.Range("C349").Value = "2:10 PM"
$ws.Range("D349").Value = "LO3837"
$ws.Range("E349").Value = "Warsaw"
$ws.Range("F349").Value = "(WAW)"
$ws.Range("G349").Value = "LOT "
$ws.Range("H349").Value = "E75S"
$ws.Range("I349").Value = "(SP-LIL)"
$ws.Range("J349").Value = "2:01 PM"
$ws.Range("L349").Value = "0 hours, -9 minutes"

# Row 350
$ws.Range("A350").Value = 349
$ws.Range("B350").Value = "Friday, Jan 13"
$ws.Range("C350").Value = "2:20 PM"
$ws.Range("D350").Value = "FR826"
$ws.Range("E350").Value = "Venice"
$ws.Range("F350").Value = "(VCE)"
$ws.Range("G350").Value = "Ryanair "
$ws.Range("H350").Value = "B738"
$ws.Range("I350").Value = "(SP-RKM)"
$ws.Range("J350").Value = "2:07 PM"
$ws.Range("L350").Value = "0 hours, -13 minutes"

# Row 351
$ws.Range("A351").Value = 350
$ws.Range("B351").Value = "Friday, Jan 13"
$ws.Range("C351").Value = "2:40 PM"
$ws.Range("D351").Value = "W61744"
$ws.Range("E351").Value = "Oslo"
$ws.Range("F351").Value = "(TRF)"
$ws.Range("G351").Value = "Wizz Air "
$ws.Range("H351").Value = "A320"
$ws.Range("I351").Value = "(HA-LYM)"
$ws.Range("J351").Value = "2:17 PM"
$ws.Range("L351").Value = "0 hours, -23 minutes"
